$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductionPlan")

# Row 2 (MAT_A / LINE_A): uncon_planned_qty (G) and con_planned_qty (H) -> 840, produced_qty (J) -> 798
$ws.Range("G2").Value = 840
$ws.Range("H2").Value = 840
$ws.Range("J2").Value = 798

# Row 3 (MAT_B / LINE_B): uncon_planned_qty (G) and con_planned_qty (H) -> 112, produced_qty (J) -> 99
$ws.Range("G3").Value = 112
$ws.Range("H3").Value = 112
$ws.Range("J3").Value = 99
